# #439 removed slug because auto-merge failed
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")
$ws.Activate()

# Remove the entire "Slug" column (column L) from the Products sheet.
# This shifts the following columns (Cash On Delivery, Product Image File,
# Quantity) left by one, and the now-unused "Slug" / "Cherry Mobile"
# shared strings are dropped automatically on save.
$ws.Columns("L").Delete()

$ws.Range("J8").Select()
